$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.146.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.308.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.13"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.20%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.85"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.02"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.09%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.668.79"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.373.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.79%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.060.06"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -11.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.94"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.85%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.29"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.63%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0690"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.015.74"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.87%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.47"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.49"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.534.59"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.87"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.82%  "
